$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 13 (existing) ---
$ws.Range("A13").Value = 23
$ws.Range("B13").Value = "31/12/2025 02:46"
$ws.Range("C13").Value = 515
$ws.Range("D13").Value = "Conhecimentos Específicos"
$ws.Range("E13").Value = "Layout e Arranjos Físicos"
$ws.Range("F13").Value = "Rever tipos de layout"

# --- Update row 14 (existing) ---
$ws.Range("A14").Value = 24
$ws.Range("B14").Value = "31/12/2025 02:53"
$ws.Range("C14").Value = 507
$ws.Range("D14").Value = "Conhecimentos Específicos"
$ws.Range("E14").Value = "Layout e Arranjos Físicos"
$ws.Range("F14").Value = "não me atentei a palavar VARIEDADE"

# --- New rows 15-27 ---
$ws.Range("A15").Value = 25
$ws.Range("B15").Value = "31/12/2025 06:12"
$ws.Range("C15").Value = 756
$ws.Range("D15").Value = "Conhecimentos Específicos"
$ws.Range("E15").Value = "Gestão da Cadeia de Suprimentos"
$ws.Range("F15").Value = "retirar os parenteses dos finais das alternativas"

$ws.Range("A16").Value = 26
$ws.Range("B16").Value = "31/12/2025 11:03"
$ws.Range("C16").Value = 1421
$ws.Range("D16").Value = "Contabilidade Gerencial"
$ws.Range("E16").Value = "Contabilidade Básica"
$ws.Range("F16").Value = "Corrigir o gabarito dessa questão"

$ws.Range("A17").Value = 27
$ws.Range("B17").Value = "31/12/2025 11:31"
$ws.Range("C17").Value = 123
$ws.Range("D17").Value = "Português"
$ws.Range("E17").Value = "Verbos Traiçoeiros"
$ws.Range("F17").Value = "arrumar essa questão, tirar os gabaritos da alternativa E"

$ws.Range("A18").Value = 28
$ws.Range("B18").Value = "31/12/2025 11:50"
$ws.Range("C18").Value = 275
$ws.Range("D18").Value = "Português"
$ws.Range("E18").Value = "Ambiguidade"
$ws.Range("F18").Value = "Tirar duvidas sobre essa questão se possível, a se sentar não seria o equivalente a sentando? se sim essa frase também é ambígua"

$ws.Range("A19").Value = 29
$ws.Range("B19").Value = "31/12/2025 11:55"
$ws.Range("C19").Value = 225
$ws.Range("D19").Value = "Português"
$ws.Range("E19").Value = "Regência Verbal"
$ws.Range("F19").Value = "O verbo “querer”, no sentido de “estimar”/“querer bem a”, é VTI e exige a preposição “a”"

$ws.Range("A20").Value = 30
$ws.Range("B20").Value = "31/12/2025 12:11"
$ws.Range("C20").Value = 222
$ws.Range("D20").Value = "Português"
$ws.Range("E20").Value = "Regência Verbal"
$ws.Range("F20").Value = "só errei essa questão porque tem dois QUEs nas alternativa E, e como nenhum veio destacado, considerei o primeiro, mas era preciso considerar o segundo para acertar a questão.`nCorrigir essa questão: marcar todos os QUEs, na E, destacar apenas o segundo"

$ws.Range("A21").Value = 31
$ws.Range("B21").Value = "31/12/2025 12:12"
$ws.Range("C21").Value = 222
$ws.Range("D21").Value = "Português"
$ws.Range("E21").Value = "Regência Verbal"
$ws.Range("F21").Value = "Como assim verbo VTI exige preposição?"

$ws.Range("A22").Value = 32
$ws.Range("B22").Value = "31/12/2025 12:25"
$ws.Range("C22").Value = 180
$ws.Range("D22").Value = "Português"
$ws.Range("E22").Value = "Travessão"
$ws.Range("F22").Value = "editar essa questão, a alternativa é foi cortada porque disparou um corte no algoritimo de leitura por causa da palavra `"comentário`""

$ws.Range("A23").Value = 33
$ws.Range("B23").Value = "31/12/2025 12:30"
$ws.Range("C23").Value = 308
$ws.Range("D23").Value = "Português"
$ws.Range("E23").Value = "Coesão"
$ws.Range("F23").Value = "Corrigir esse texto, quebras de linhas"

$ws.Range("A24").Value = 34
$ws.Range("B24").Value = "31/12/2025 12:32"
$ws.Range("C24").Value = 294
$ws.Range("D24").Value = "Português"
$ws.Range("E24").Value = "Coesão"
$ws.Range("F24").Value = "Arrumar esse tal destaque. Deve ser no pronome `"lo`""

$ws.Range("A25").Value = 35
$ws.Range("B25").Value = "31/12/2025 12:35"
$ws.Range("C25").Value = 288
$ws.Range("D25").Value = "Português"
$ws.Range("E25").Value = "Coesão"
$ws.Range("F25").Value = "extrair o texto dessa imagem e vinculá-lo a essa questão"

$ws.Range("A26").Value = 36
$ws.Range("B26").Value = "31/12/2025 12:42"
$ws.Range("C26").Value = 281
$ws.Range("D26").Value = "Português"
$ws.Range("E26").Value = "Coesão"
$ws.Range("F26").Value = "Tentar identificar essas tais palavras em destaque, acertei essa questão prativamente dando um tiro no escuro"

$ws.Range("A27").Value = 37
$ws.Range("B27").Value = "31/12/2025 12:45"
$ws.Range("C27").Value = 303
$ws.Range("D27").Value = "Português"
$ws.Range("E27").Value = "Reescritura"
$ws.Range("F27").Value = "Extrair o texto dessa imagem e vinculá-lo à questão"

Write-Output "done"
